$d = $word.ActiveDocument
$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>
    <w:p>
      <w:pPr>
        <w:jc w:val="right"/>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
        <w:t xml:space="preserve">CONVOCATORIA – </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
        <w:t>{carrera}</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
        <w:t xml:space="preserve">– </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
        <w:t>{</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
        <w:t>fechaActual</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
        <w:t>}</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:bookmarkStart w:id="0" w:name="_GoBack"/>
      <w:bookmarkEnd w:id="0"/>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="right"/>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
        <w:t xml:space="preserve">Cuenca, </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
        <w:t>{</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
        <w:t>fechaActual</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
        <w:t>}</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
        <w:t>A los estudiantes Interesados:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
        <w:t>Se convoca a los estudiantes de</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
        <w:t xml:space="preserve"> {curso}</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
        <w:t>. ciclo de la carrera</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
        <w:t xml:space="preserve"> {carrera}</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
        <w:t xml:space="preserve"> que deseen realizar sus prácticas pre profesionales en la empresa</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
        <w:t xml:space="preserve"> {empresa}</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
        <w:t>, a presentar la solicitud correspondiente.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
        <w:t>Las actividades a desarrollar son:</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
        <w:t xml:space="preserve"> desarrollar software para móviles</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
        <w:t>, por lo que los postulantes deberán haber aprobado las siguientes asignaturas:</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
        <w:t xml:space="preserve"> {materia}</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
        <w:t>La fecha máxima en la que se receptarán las solicitudes es el</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
        <w:t xml:space="preserve"> {</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
        <w:t>fechaLimite</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
        <w:t>}</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
        <w:t xml:space="preserve">Para mayor información contactarse con </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
        <w:t>la secretaría del Instituto</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
        <w:t>, docente responsable de prácticas pre profesionales de la carrera.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
        <w:t>Atentamente,</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
        <w:t>______________________</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
        <w:t>Responsable de Prácticas Pre Profesionales</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
        <w:t xml:space="preserve">CARRERA DE </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
        <w:t>{carrera}</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:bCs/>
        </w:rPr>
        <w:t>INSTITUTO SUPERIOR TECNOLÓGICO DEL AZUAY</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:sectPr>
      <w:pgSz w:w="11906" w:h="16838"/>
      <w:pgMar w:top="1417" w:right="1701" w:bottom="1417" w:left="1701" w:header="708" w:footer="708" w:gutter="0"/>
      <w:cols w:space="708"/>
      <w:docGrid w:linePitch="360"/>
    </w:sectPr>
  </w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Content.InsertXML($xml)
Write-Output ("Paragraphs: " + $d.Paragraphs.Count)
